$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "25.874.65";     E = "  -2.31%  " },
    @{ Row = 3;  D = "1.754.25";      E = "  -4.69%  " },
    @{ Row = 4;  D = "1.002";         E = "  +0.17%  " },
    @{ Row = 5;  D = "239.34";        E = "  -8.27%  " },
    @{ Row = 6;  D = $null;           E = "  +0.07%  " },
    @{ Row = 7;  D = "0.5078";        E = "  -5.64%  " },
    @{ Row = 8;  D = "42.44";         E = "  -5.27%  " },
    @{ Row = 9;  D = "0.2764";        E = "  -6.36%  " },
    @{ Row = 10; D = "0.06190";       E = "  -10.97%  " },
    @{ Row = 11; D = "1.751.80";      E = "  -4.81%  " },
    @{ Row = 12; D = "0.06972";       E = "  -3.30%  " },
    @{ Row = 13; D = "15.65";         E = "  -10.05%  " },
    @{ Row = 14; D = "0.6128";        E = "  -15.90%  " },
    @{ Row = 15; D = "4.529";         E = "  -9.09%  " },
    @{ Row = 16; D = "77.31";         E = "  -13.39%  " },
    @{ Row = 17; D = $null;           E = "  +0.13%  " },
    @{ Row = 18; D = "1.001";         E = "  +0.07%  " },
    @{ Row = 19; D = "25.881.92";     E = "  -2.34%  " },
    @{ Row = 20; D = "0.000006869";   E = "  -12.98%  " },
    @{ Row = 21; D = "11.70";         E = "  -15.05%  " },
    @{ Row = 22; D = "1.975.38";      E = "  -4.92%  " },
    @{ Row = 23; D = "4.076";         E = "  -11.14%  " },
    @{ Row = 24; D = "8.251";         E = "  -10.07%  " },
    @{ Row = 25; D = "5.243";         E = "  -12.45%  " },
    @{ Row = 26; D = "137.93";        E = "  -3.64%  " },
    @{ Row = 27; D = "1.489";         E = "  -12.86%  " },
    @{ Row = 28; D = $null;           E = "  -11.16%  " },
    @{ Row = 29; D = "1.817";         E = "  -15.85%  " },
    @{ Row = 30; D = "103.68";        E = "  -6.64%  " },
    @{ Row = 31; D = "0.08238";       E = "  -7.47%  " },
    @{ Row = 32; D = $null;           E = "  -13.30%  " },
    @{ Row = 33; D = "3.491";         E = "  -13.50%  " },
    @{ Row = 34; D = "0.04554";       E = "  -5.94%  " },
    @{ Row = 35; D = "1.000";         E = "  +0.14%  " },
    @{ Row = 36; D = "2.644";         E = "  -8.81%  " },
    @{ Row = 37; D = "0.9947";        E = "  -12.12%  " },
    @{ Row = 38; D = "0.6099";        E = "  -16.09%  " },
    @{ Row = 39; D = "2.703";         E = "  -12.71%  " },
    @{ Row = 40; D = $null;           E = "  -8.71%  " },
    @{ Row = 41; D = $null;           E = "  +0.12%  " },
    @{ Row = 42; D = "1.899";         E = "  -17.41%  " },
    @{ Row = 43; D = "103.34";        E = "  -3.78%  " },
    @{ Row = 44; D = "0.3854";        E = "  -17.65%  " },
    @{ Row = 45; D = "0.7403";        E = "  -18.14%  " },
    @{ Row = 46; D = "4.983";         E = "  -15.30%  " },
    @{ Row = 47; D = "0.05440";       E = "  -5.17%  " },
    @{ Row = 48; D = "0.1113";        E = "  -10.37%  " },
    @{ Row = 49; D = "6.017";         E = "  -19.09%  " },
    @{ Row = 50; D = "7.700";         E = "  -14.38%  " },
    @{ Row = 51; D = "30.05";         E = "  -13.68%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($u.Row, 4)
        $origStyleD = $cellD.Style
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.Style = $origStyleD
    }
    $cellE = $ws.Cells.Item($u.Row, 5)
    $origStyleE = $cellE.Style
    $cellE.NumberFormat = "@"
    $cellE.Value = $u.E
    $cellE.Style = $origStyleE
}
